$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values. Some new values are plain numeric
# strings (e.g. "568.99", "1.00"); without forcing a Text number format
# first, Excel would parse them as numbers and normalise/trim the text
# (e.g. "134.00" -> 134), which would not match the source data which
# stores these as literal text strings.
$ws.Range("D2").Value = "61.419.78"
$ws.Range("D3").Value = "2.890.07"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.99"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.84"
$ws.Range("D9").Value = "2.887.40"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("D16").Value = "3.369.34"
$ws.Range("D17").Value = "61.473.23"
$ws.Range("D19").Value = "2.881.20"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.32"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.38"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.79"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000103"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.72"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.24"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.66"
$ws.Range("D45").Value = "2.684.26"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.00"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "339.49"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.51"

# Update "Volume(1h)" column (E) values (padded percentage strings).
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -4.59%  "
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("E10").Value = "  -8.44%  "
$ws.Range("E11").Value = "  -4.92%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -11.90%  "
$ws.Range("E28").Value = "  -7.88%  "
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E31").Value = "  -5.37%  "
$ws.Range("E32").Value = "  -8.23%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E39").Value = "  -13.38%  "
$ws.Range("E40").Value = "  -6.71%  "
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -8.39%  "
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  -7.51%  "
